# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the regenerated output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map of row -> new F value, identical update set applied to both sheets
# (the "全部类型" sheet has one extra row inserted earlier in the data,
# so its row numbers differ by +1 from row 33 onward).
$updatesSheet1 = @{
    4  = 1137
    6  = 83
    8  = 60
    9  = 1170
    10 = 16452
    12 = 205
    13 = 1042
    14 = 6388
    20 = 1276
    23 = 639
    28 = 227
    29 = 899
    31 = 5062
    32 = 506
    33 = 11372
    34 = 1251
    35 = 20
    36 = 155
    37 = 213
    38 = 3845
    39 = 272
}

$updatesSheet4 = @{
    4  = 1137
    6  = 83
    8  = 60
    9  = 1170
    10 = 16452
    12 = 205
    13 = 1042
    14 = 6388
    20 = 1276
    23 = 639
    28 = 227
    29 = 899
    31 = 5062
    32 = 506
    34 = 11372
    35 = 1251
    36 = 20
    37 = 155
    38 = 213
    39 = 3845
    40 = 272
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updatesSheet1.Keys) {
    $ws1.Range("F$row").Value = $updatesSheet1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesSheet4.Keys) {
    $ws4.Range("F$row").Value = $updatesSheet4[$row]
}
